# Finished macsj2243 and macsj2228 reductions:
# mark the Flat Reduced / Image / SExtractor Catalog columns ("wd") as
# complete for several rows, add a new comment about the deep band used,
# and move the active selection to reflect where work left off.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("MC3PO Sample")

# Row 12: Flat Reduced (K), Image (L), SExtractor Catalog (M)
$ws.Range("K12").Value = "wd"
$ws.Range("L12").Value = "wd"
$ws.Range("M12").Value = "wd"

# Row 13: Flat Reduced (K), Image (L), SExtractor Catalog (M)
$ws.Range("K13").Value = "wd"
$ws.Range("L13").Value = "wd"
$ws.Range("M13").Value = "wd"

# Row 14: Image (L), SExtractor Catalog (M)
$ws.Range("L14").Value = "wd"
$ws.Range("M14").Value = "wd"

# Row 15: Image (L), SExtractor Catalog (M), plus a new comment in P15
$ws.Range("L15").Value = "wd"
$ws.Range("M15").Value = "wd"
$ws.Range("P15").Value = "I used this as the deepband during reduction since the 2000aug RC band only had 8 chips."

# Row 16: Image (L), SExtractor Catalog (M)
$ws.Range("L16").Value = "wd"
$ws.Range("M16").Value = "wd"

# Row 17: Image (L), SExtractor Catalog (M)
$ws.Range("L17").Value = "wd"
$ws.Range("M17").Value = "wd"

# Leave the selection where the author's work left off
$ws.Range("M11").Select()
